$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update nombre_aides (column C) and montant_total (column D) for the 2020-08-19 data refresh
$updates = @(
    @{Row=2; C=38024; D=54989645}
    @{Row=3; C=91539; D=134183736}
    @{Row=4; C=31333; D=46404154}
    @{Row=5; C=8750; D=13005563}
    @{Row=6; C=2012; D=2989971}
    @{Row=12; C=41553; D=56380616}
    @{Row=13; C=9734; D=14080636}
    @{Row=14; C=26091; D=38260073}
    @{Row=15; C=8344; D=12383824}
    @{Row=16; C=2160; D=3212165}
    @{Row=17; C=420; D=619123}
    @{Row=20; C=10282; D=13607382}
    @{Row=21; C=13492; D=19475595}
    @{Row=22; C=31815; D=46687227}
    @{Row=23; C=10263; D=15256078}
    @{Row=24; C=2657; D=3950682}
    @{Row=25; C=510; D=759092}
    @{Row=27; C=11743; D=15684114}
    @{Row=28; C=7721; D=11179237}
    @{Row=29; C=22634; D=33226036}
    @{Row=30; C=7861; D=11699133}
    @{Row=31; C=1981; D=2955999}
    @{Row=34; C=8360; D=11041631}
    @{Row=35; C=3272; D=4723937}
    @{Row=36; C=7897; D=11532912}
    @{Row=37; C=3193; D=4732461}
    @{Row=41; C=2490; D=3366932}
    @{Row=42; C=17378; D=25127402}
    @{Row=43; C=51435; D=75404654}
    @{Row=44; C=19097; D=28366508}
    @{Row=45; C=5638; D=8393305}
    @{Row=46; C=1214; D=1811545}
    @{Row=47; C=63; D=92568}
    @{Row=50; C=16805; D=22366649}
    @{Row=51; C=2054; D=2979454}
    @{Row=52; C=6984; D=10265599}
    @{Row=53; C=2373; D=3543964}
    @{Row=54; C=758; D=1132305}
    @{Row=55; C=187; D=277333}
    @{Row=57; C=7089; D=9744112}
    @{Row=58; C=1014; D=1573066}
    @{Row=59; C=2541; D=3976151}
    @{Row=60; C=1014; D=1602039}
    @{Row=61; C=344; D=543883}
    @{Row=62; C=112; D=179850}
    @{Row=63; C=20; D=33000}
    @{Row=64; C=1468; D=2142823}
    @{Row=65; C=15495; D=22380327}
    @{Row=66; C=45002; D=65848295}
    @{Row=67; C=15777; D=23445042}
    @{Row=68; C=4587; D=6831551}
    @{Row=69; C=931; D=1384668}
    @{Row=73; C=15182; D=20005836}
    @{Row=74; C=52367; D=76214967}
    @{Row=75; C=147947; D=217973316}
    @{Row=76; C=64129; D=95561352}
    @{Row=77; C=20504; D=30636322}
    @{Row=78; C=4873; D=7278543}
    @{Row=79; C=269; D=398670}
    @{Row=85; C=51612; D=70213695}
    @{Row=86; C=4650; D=6738489}
    @{Row=87; C=11661; D=17131527}
    @{Row=88; C=3904; D=5818583}
    @{Row=89; C=1351; D=2018989}
    @{Row=93; C=5447; D=7322009}
    @{Row=94; C=1609; D=2317033}
    @{Row=95; C=5215; D=7682135}
    @{Row=96; C=1945; D=2897437}
    @{Row=97; C=693; D=1038460}
    @{Row=98; C=187; D=279613}
    @{Row=101; C=3590; D=4751335}
    @{Row=102; C=643; D=1007439}
    @{Row=103; C=381; D=608380}
    @{Row=104; C=136; D=209660}
    @{Row=105; C=47; D=73500}
    @{Row=106; C=22; D=36000}
    @{Row=107; C=10858; D=15753660}
    @{Row=108; C=29382; D=43163618}
    @{Row=109; C=9830; D=14617679}
    @{Row=110; C=2706; D=4035207}
    @{Row=111; C=494; D=736046}
    @{Row=113; C=7; D=10500}
    @{Row=114; C=9843; D=13001712}
    @{Row=115; C=30697; D=44267742}
    @{Row=116; C=66511; D=97337095}
    @{Row=117; C=21466; D=31902419}
    @{Row=118; C=6090; D=9073021}
    @{Row=119; C=1133; D=1693271}
    @{Row=120; C=79; D=115920}
    @{Row=124; C=25985; D=34703726}
    @{Row=125; C=36279; D=52361014}
    @{Row=126; C=77213; D=112907353}
    @{Row=127; C=23965; D=35568356}
    @{Row=128; C=6420; D=9540738}
    @{Row=129; C=1246; D=1853411}
    @{Row=133; C=31963; D=42438982}
    @{Row=134; C=13346; D=19319381}
    @{Row=135; C=32500; D=47732929}
    @{Row=136; C=11522; D=17119892}
    @{Row=137; C=2972; D=4430714}
    @{Row=138; C=502; D=746990}
    @{Row=141; C=10870; D=14491983}
    @{Row=142; C=35357; D=51063454}
    @{Row=143; C=81808; D=119855010}
    @{Row=144; C=24475; D=36362348}
    @{Row=145; C=6429; D=9593067}
    @{Row=146; C=1448; D=2154230}
    @{Row=147; C=82; D=122630}
    @{Row=149; C=29354; D=39592308}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    $ws.Cells.Item($u.Row, 4).Value = $u.D
}
